# Word COM-interop script: minor formatting changes to the tutorial document.
#
# 1. Insert a new, completely empty paragraph right after the paragraph that
#    ends in " like this:" (and right before the centered paragraph that
#    holds the screenshot which follows it).
#
# 2. Insert a new, empty paragraph (carrying "spacing after=0 / justify
#    both" formatting) right before the existing blank paragraph that
#    precedes the centered screenshot paragraph that follows the text
#    "... of the dialog:".

$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("like this:", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Text = "`r"

# --- Change 2 --------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("of the dialog", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$rng2.Collapse(0)
$rng2.MoveEnd(1, 1) | Out-Null
$rng2.Collapse(0)
$rng2.Text = "`r"
